$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new column to the left of the existing column A.
#    This shifts: A->B (snowdepth values), B->C (density mean values),
#    C->D (density SLF values), D->E (density1 values), E->F, F->G
$ws.Columns.Item(1).Insert()

# 2) Insert a new row above the existing row 2 (data start), shifting
#    all existing data rows down by one.
$ws.Rows.Item(2).Insert()

# 3) Header row (row 1)
$ws.Range("A1").Value = "snowdepth"
$ws.Range("B1").Value = "snowheight"
$ws.Range("C1").Value = "mean"
$ws.Range("D1").Value = "density SLF in kg/m^3"
$ws.Range("E1").Value = "density1 "

# 4) New snowdepth column (A) values for every measurement row
$snowdepth = @(0, 8, 18, 28, 38, 48, 58, 68, 78, 88, 98)
for ($i = 0; $i -lt $snowdepth.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 1).Value = $snowdepth[$i]
}

# 5) New row 2 data (depth=0 measurement) in columns B:D
$ws.Range("B2").Value = 118
$ws.Range("C2").Value = 100
$ws.Range("D2").Value = 100

# 6) Update the active selection to C1, matching the saved view state
$ws.Range("C1").Select()
